# Trade #15 closed at 2026-02-17 12:28:15 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" aggregate metrics to account
# for the newly-closed MarketMaking trade, and appends the trade's row to
# both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.01   # Current Capital
$summary.Range("B4").Value = 0.01      # Total P&L $
$summary.Range("B5").Value = 0.01      # Total P&L %
$summary.Range("B6").Value = 15        # Total Trades
$summary.Range("B8").Value = 6         # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.01     # Capital
$status.Range("D4").Value = 15         # Trades
$status.Range("E4").Value = 0.01       # P&L $
$status.Range("F4").Value = 0.01       # P&L %
$status.Range("G4").Value = 40         # Win Rate %

# ---------------------------------------------------------------------
# Helper: write a text value into a cell without Excel re-interpreting
# strings that look like dates/times as date serials.
# ---------------------------------------------------------------------
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Append the new closed trade (row 16) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(16, 1).Value = 15                 # Trade #
    Set-TextCell $ws 16 2 "2026-02-17"                # Date
    Set-TextCell $ws 16 3 "12:28:09"                  # Time
    Set-TextCell $ws 16 4 "MarketMaking"               # Strategy
    Set-TextCell $ws 16 5 "UP"                         # Side
    $ws.Cells.Item(16, 6).Value = 0.08                # Entry Price
    $ws.Cells.Item(16, 7).Value = 0.06                # Exit Price
    Set-TextCell $ws 16 8 "CLOSED"                     # Status
    $ws.Cells.Item(16, 9).Value = -25                 # P&L %
    $ws.Cells.Item(16, 10).Value = -0.02              # P&L $
    $ws.Cells.Item(16, 11).Value = 100.01             # Capital After
    $ws.Cells.Item(16, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(16, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(16, 14).Value = 0.6                # Confidence
    Set-TextCell $ws 16 15 "Normal spread capture: 19600 bps"  # Entry Reason
    Set-TextCell $ws 16 16 "early_exit"                # Exit Reason
    $ws.Cells.Item(16, 17).Value = 0.14               # Duration (min)
}
